$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mIF")
$ws.Rows.Item(4).Delete()
